$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Investment_Cost")

# --- New "Lifetime" column header (same style as the other header cells) ---
$ws.Range("G1").Value = "Lifetime"
$ws.Range("G1").Font.Bold = $true
$ws.Range("G1").WrapText = $true

# --- Per-row lifetimes (column G) ---
$ws.Range("G2").Value = "35Y"    # PV_plant
$ws.Range("G3").Value = "25Y"    # AEC_Electrolyzer
$ws.Range("G4").Value = "25Y"    # PEM_Electrolyzer
$ws.Range("G5").Value = "25Y"    # SOEC_Electrolyzer
$ws.Range("G6").Value = "20Y"    # CO2_Vaporizer
$ws.Range("G7").Value = "30Y"    # Methanol_Plant
$ws.Range("G8").Value = "20Y"    # Electric_Steam_Boiler
$ws.Range("G9").Value = "25Y"    # Methanol_storage
$ws.Range("G10").Value = "25Y"   # Hydrogen_storage
$ws.Range("G11").Value = "27Y"   # Wind_onshore (was Wind_farm)
$ws.Range("G12").Value = "27Y"   # Wind_offshore (new)
$ws.Range("G13").Value = "30Y"   # Destilation_tower (new)

# --- Rename Wind_farm -> Wind_onshore, update its costs ---
$ws.Range("A11").Value = "Wind_onshore"
$ws.Range("B11").Value = 1110000
$ws.Range("C11").Value = 1180000
$ws.Range("D11").Value = 1150000
$ws.Range("E11").Value = 1110000
$ws.Range("F11").Value = 1090000

# --- New row: Wind_offshore ---
$ws.Range("A12").Value = "Wind_offshore"
$ws.Range("B12").Value = 2120000
$ws.Range("C12").Value = 1880000
$ws.Range("D12").Value = 1800000
$ws.Range("E12").Value = 1680000
$ws.Range("F12").Value = 1640000

# --- New row: Destilation_tower ---
$ws.Range("A13").Value = "Destilation_tower"
$ws.Range("B13").Value = 1350000
$ws.Range("C13").Value = 1350000
$ws.Range("D13").Value = 1090000
$ws.Range("E13").Value = 960000
$ws.Range("F13").Value = 870000

# --- Apply the new number format to the investment-cost cells ---
# (rows 9 and 10 keep their original 164 number format/style)
$ws.Range("B2:F8").NumberFormat = "#,##0.00"
$ws.Range("B11:F13").NumberFormat = "#,##0.00"

# --- Make Investment_Cost the active sheet / selection ---
$ws.Activate()
$ws.Range("G5").Select()
